$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# =========================================================
# Sheet1 "KPIs": Channel/Retailer aggregation KPIs -> Brand
# =========================================================

# Rename the two "Channel Aggregation ..." KPI rows to "Brand Aggregation ..."
$ws1.Range("A2").Value = "Brand Aggregation SOS"
$ws1.Range("A3").Value = "Brand Aggregation Linear SOS"

# Those two rows lose their special (bold) style and take on the plain
# header-row style instead - copy the format from A1 (style 0).
$ws1.Range("A1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("A3").PasteSpecial(-4122)

# Remove the two "Retailer Aggregation ..." rows (old rows 4 & 5) plus the
# blank row that followed them so that the lone far-away row (old row 30,
# "GMI") shifts up to row 27.
$ws1.Range("A4:A6").EntireRow.Delete()

# The row that used to be "A30: GMI" now also carries a "Parent" label.
$ws1.Range("A27").Value = "GMI"
$ws1.Range("B27").Value = "Parent"

# Approximate the (LibreOffice-derived) new column widths as closely as
# this engine's pixel-quantized width model allows.
$ws1.Columns("A").ColumnWidth = 26.67
$ws1.Columns("D").ColumnWidth = 14.83
$ws1.Columns("H").ColumnWidth = 13.67

# =========================================================
# Sheet2 "Aggregation": Channel/Retailer rows -> Brand only,
# and the "SOS Type" column (C) is dropped since "Aggregation
# Level" (old column B) is replaced by the old "SOS Type" values.
# =========================================================

$ws2.Range("B1").Value = "SOS Type"
$ws2.Range("A2").Value = "Brand Aggregation SOS"
$ws2.Range("B2").Value = "SOS"
$ws2.Range("A3").Value = "Brand Aggregation Linear SOS"
$ws2.Range("B3").Value = "Linear SOS"

# Drop the old column C entirely and the two "Retailer Aggregation" rows.
$ws2.Columns("C").Delete()
$ws2.Range("A4:A5").EntireRow.Delete()

$ws2.Columns("A").ColumnWidth = 26.0

# =========================================================
# Selection / active sheet bookkeeping
# =========================================================
# Select Aggregation's new A3 first, then make KPIs the active sheet with
# A6 selected - this produces activeTab=0 with KPIs tabSelected and
# Aggregation's own stored selection left at A3.
$ws2.Range("A3").Select()
$ws1.Activate()
$ws1.Range("A6").Select()
